$p = $ppt.ActivePresentation

# --- Slide 4: "informações" -> "notícias" (first paragraph of content placeholder) ---
$s4 = $p.Slides.Item(4)
$tr4 = $s4.Shapes.Item(2).TextFrame.TextRange
$para4 = $tr4.Paragraphs(1, 1)
$para4.Text = "Como um cliente de corretoras de criptomoedas, eu quero ver notícias sobre as moedas disponíveis para mineração, para que eu possa investir "

# --- Slide 5: update the user story wording ---
$s5 = $p.Slides.Item(5)
$tf5 = $s5.Shapes.Item(2).TextFrame
$tf5.TextRange.Text = "Eu enquanto minerador, gostaria de uma solução que monitorasse minha máquina para avaliar o desempenho da minha GPU ao longo da mineração "
$tf5.TextRange.LanguageID = "pt-BR"

# --- New slide 7: add a new "Título e Conteúdo" slide at the end ---
$s7 = $p.Slides.Add(7, 16)

$tf1 = $s7.Shapes.Item(1).TextFrame
$tf1.TextRange.Text = "4"
$tf1.TextRange.LanguageID = "pt-BR"

$tf2 = $s7.Shapes.Item(2).TextFrame
$tf2.TextRange.Text = "Eu enquanto minerador iniciante, gostaria de ter acesso a um teste prévio de uma ferramenta de monitoramento para ter maior confiança para compra-lá "
$tf2.TextRange.LanguageID = "pt-BR"
